$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Type Of Inventory" value for the first data row first
$ws.Range("G2").Value = "Raw Material"

# Fill in example values across the rest of the first data row
$ws.Range("A2:F2").Value = "example"

# Extend the existing date formatting down into B6 (same style as B2:B5)
$ws.Range("B2").Copy()
$ws.Range("B6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Leave the selection where the author last clicked
$ws.Range("C13").Select()
